$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 10.18201399862616
$ws.Range("C2").Value = 4.340856429089271
$ws.Range("D2").Value = 8.631474025476477
$ws.Range("E2").Value = 13.63263849530765
$ws.Range("F2").Value = 36.21285343707838
$ws.Range("J2").Value = 10.32451255953204
$ws.Range("K2").Value = 9.564593265398603
$ws.Range("M2").Value = 15.02640002581191
$ws.Range("O2").Value = 27.63659263531571

# Row 3
$ws.Range("B3").Value = 9.943804821937535
$ws.Range("C3").Value = 4.194357348631216
$ws.Range("D3").Value = 8.599908211621893
$ws.Range("E3").Value = 13.62914680221801
$ws.Range("F3").Value = 36.28459783135274
$ws.Range("J3").Value = 10.34531993874193
$ws.Range("K3").Value = 9.407726021104814
$ws.Range("M3").Value = 14.96942325103585
$ws.Range("O3").Value = 27.71900279472649

# Row 4
$ws.Range("B4").Value = 9.796481557039414
$ws.Range("C4").Value = 4.101062199758176
$ws.Range("D4").Value = 8.581810633105528
$ws.Range("E4").Value = 13.62927148237988
$ws.Range("F4").Value = 36.33633611243493
$ws.Range("J4").Value = 10.35926522512753
$ws.Range("K4").Value = 9.311584336855518
$ws.Range("M4").Value = 14.93664562874724
$ws.Range("O4").Value = 27.77486253062056

# Row 5
$ws.Range("B5").Value = 9.736269120299138
$ws.Range("C5").Value = 4.062240146958278
$ws.Range("D5").Value = 8.574763920546062
$ws.Range("E5").Value = 13.62989428276155
$ws.Range("F5").Value = 36.35934970794574
$ws.Range("J5").Value = 10.36524240525651
$ws.Range("K5").Value = 9.272498714816226
$ws.Range("M5").Value = 14.92385294399519
$ws.Range("O5").Value = 27.7989464541241

# Row 6
$ws.Range("B6").Value = 9.726262896112365
$ws.Range("C6").Value = 4.055746414388171
$ws.Range("D6").Value = 8.573613790983462
$ws.Range("E6").Value = 13.63003228520333
$ws.Range("F6").Value = 36.36328755531196
$ws.Range("J6").Value = 10.3662526986388
$ws.Range("K6").Value = 9.266015621596335
$ws.Range("M6").Value = 14.92176309133655
$ws.Range("O6").Value = 27.80302527584623

# Row 7
$ws.Range("B7").Value = 9.795670106358765
$ws.Range("C7").Value = 4.100541833099856
$ws.Range("D7").Value = 8.581714263008898
$ws.Range("E7").Value = 13.62927756386048
$ws.Range("F7").Value = 36.33663867249513
$ws.Range("J7").Value = 10.3593446432253
$ws.Range("K7").Value = 9.311056773591565
$ws.Range("M7").Value = 14.93647080431495
$ws.Range("O7").Value = 27.77518199029936

# Row 8
$ws.Range("B8").Value = 10.10015277305855
$ws.Range("C8").Value = 4.291058831482843
$ws.Range("D8").Value = 8.620326945907751
$ws.Range("E8").Value = 13.63096470493369
$ws.Range("F8").Value = 36.23599372932272
$ws.Range("J8").Value = 10.33144439581669
$ws.Range("K8").Value = 9.510498130960714
$ws.Range("M8").Value = 15.00630247138895
$ws.Range("O8").Value = 27.66391489844211

# Row 9
$ws.Range("B9").Value = 10.68505201329504
$ws.Range("C9").Value = 4.636682066171169
$ws.Range("D9").Value = 8.705985580188742
$ws.Range("E9").Value = 13.65219697096628
$ws.Range("F9").Value = 36.09974565224459
$ws.Range("J9").Value = 10.28599846796049
$ws.Range("K9").Value = 9.900883870613859
$ws.Range("M9").Value = 15.16029985778427
$ws.Range("O9").Value = 27.48753658605569

# Row 10
$ws.Range("B10").Value = 11.10260869652087
$ws.Range("C10").Value = 4.871874571669933
$ws.Range("D10").Value = 8.77462622582383
$ws.Range("E10").Value = 13.67860724859336
$ws.Range("F10").Value = 36.03704875132539
$ws.Range("J10").Value = 10.25824086807658
$ws.Range("K10").Value = 10.1845072421317
$ws.Range("M10").Value = 15.28320591447971
$ws.Range("O10").Value = 27.38355434984816

# Row 11
$ws.Range("B11").Value = 11.28901741444837
$ws.Range("C11").Value = 4.974512900817066
$ws.Range("D11").Value = 8.807010973570232
$ws.Range("E11").Value = 13.69293948522034
$ws.Range("F11").Value = 36.01667000668015
$ws.Range("J11").Value = 10.24683223491522
$ws.Range("K11").Value = 10.31228527292682
$ws.Range("M11").Value = 15.34109378105421
$ws.Range("O11").Value = 27.34183146054809

# Row 12
$ws.Range("B12").Value = 11.35902794234373
$ws.Range("C12").Value = 5.012733402440843
$ws.Range("D12").Value = 8.819433755927628
$ws.Range("E12").Value = 13.69869705613012
$ws.Range("F12").Value = 36.01012481943727
$ws.Range("J12").Value = 10.24268697846641
$ws.Range("K12").Value = 10.36044858737724
$ws.Range("M12").Value = 15.36328573881582
$ws.Range("O12").Value = 27.32683601417247

# Row 13
$ws.Range("B13").Value = 11.34397669593823
$ws.Range("C13").Value = 5.004530982569888
$ws.Range("D13").Value = 8.816751319081808
$ws.Range("E13").Value = 13.69744242336718
$ws.Range("F13").Value = 36.01148231738225
$ws.Range("J13").Value = 10.24357195716832
$ws.Range("K13").Value = 10.35008638274629
$ws.Range("M13").Value = 15.35849446245176
$ws.Range("O13").Value = 27.33002976113037

# Row 14
$ws.Range("B14").Value = 11.29478919350358
$ws.Range("C14").Value = 4.977670388009038
$ws.Range("D14").Value = 8.808029852268758
$ws.Range("E14").Value = 13.69340656359716
$ws.Range("F14").Value = 36.01610804213654
$ws.Range("J14").Value = 10.24648769755093
$ws.Range("K14").Value = 10.31625242446066
$ws.Range("M14").Value = 15.34291417123239
$ws.Range("O14").Value = 27.34058165060211

# Row 15
$ws.Range("B15").Value = 11.26458313325095
$ws.Range("C15").Value = 4.961132751088233
$ws.Range("D15").Value = 8.802708234843042
$ws.Range("E15").Value = 13.69097739738849
$ws.Range("F15").Value = 36.01909405363057
$ws.Range("J15").Value = 10.24829644778471
$ws.Range("K15").Value = 10.29549775832896
$ws.Range("M15").Value = 15.33340569552185
$ws.Range("O15").Value = 27.34714975825482

# Row 16
$ws.Range("B16").Value = 11.09034916586434
$ws.Range("C16").Value = 4.865077520724664
$ws.Range("D16").Value = 8.772532517368029
$ws.Range("E16").Value = 13.67771699220092
$ws.Range("F16").Value = 36.03854448831284
$ws.Range("J16").Value = 10.25901094473569
$ws.Range("K16").Value = 10.17612776209553
$ws.Range("M16").Value = 15.27946147051785
$ws.Range("O16").Value = 27.38639346482496

# Row 17
$ws.Range("B17").Value = 10.98250458335518
$ws.Range("C17").Value = 4.80502059141437
$ws.Range("D17").Value = 8.754312552866006
$ws.Range("E17").Value = 13.67017377527321
$ws.Range("F17").Value = 36.05256299944557
$ws.Range("J17").Value = 10.26589582702949
$ws.Range("K17").Value = 10.10254723044531
$ws.Range("M17").Value = 15.24686556467192
$ws.Range("O17").Value = 27.41189850271012

# Row 18
$ws.Range("B18").Value = 10.92014585289805
$ws.Range("C18").Value = 4.770068861941884
$ws.Range("D18").Value = 8.743942532575398
$ws.Range("E18").Value = 13.66605350478012
$ws.Range("F18").Value = 36.06139245471574
$ws.Range("J18").Value = 10.2699705307049
$ws.Range("K18").Value = 10.06011114236035
$ws.Range("M18").Value = 15.22830419283106
$ws.Range("O18").Value = 27.42709335984739

# Row 19
$ws.Range("B19").Value = 10.89897788718379
$ws.Range("C19").Value = 4.758165306472882
$ws.Range("D19").Value = 8.740450470875253
$ws.Range("E19").Value = 13.66469604962916
$ws.Range("F19").Value = 36.06451354542017
$ws.Range("J19").Value = 10.27136986338277
$ws.Range("K19").Value = 10.0457247365301
$ws.Range("M19").Value = 15.22205213408523
$ws.Range("O19").Value = 27.43232819476685

# Row 20
$ws.Range("B20").Value = 10.99401943199561
$ws.Range("C20").Value = 4.811456197522268
$ws.Range("D20").Value = 8.7562408066103
$ws.Range("E20").Value = 13.67095418262649
$ws.Range("F20").Value = 36.05099137918315
$ws.Range("J20").Value = 10.2651510502747
$ws.Range("K20").Value = 10.1103921848074
$ws.Range("M20").Value = 15.25031620257183
$ws.Range("O20").Value = 27.40912909518175

# Row 21
$ws.Range("B21").Value = 11.3092529798633
$ws.Range("C21").Value = 4.985577690055722
$ws.Range("D21").Value = 8.810587292969405
$ws.Range("E21").Value = 13.69458305596565
$ws.Range("F21").Value = 36.01471754756044
$ws.Range("J21").Value = 10.24562652810364
$ws.Range("K21").Value = 10.32619668810455
$ws.Range("M21").Value = 15.34748323096325
$ws.Range("O21").Value = 27.33746046767084

# Row 22
$ws.Range("B22").Value = 11.51187645310357
$ws.Range("C22").Value = 5.095600668217984
$ws.Range("D22").Value = 8.847031274620463
$ws.Range("E22").Value = 13.71194931071946
$ws.Range("F22").Value = 35.9978407002079
$ws.Range("J22").Value = 10.23388571106605
$ws.Range("K22").Value = 10.46591571661023
$ws.Range("M22").Value = 15.41256114609167
$ws.Range("O22").Value = 27.29530858694986

# Row 23
$ws.Range("B23").Value = 11.40406513968337
$ws.Range("C23").Value = 5.037230922739933
$ws.Range("D23").Value = 8.827498260510138
$ws.Range("E23").Value = 13.70250570362727
$ws.Range("F23").Value = 36.00622306062057
$ws.Range("J23").Value = 10.24005880328268
$ws.Range("K23").Value = 10.39148018565443
$ws.Range("M23").Value = 15.37768832538335
$ws.Range("O23").Value = 27.31737633956313

# Row 24
$ws.Range("B24").Value = 10.98881467695224
$ws.Range("C24").Value = 4.808547979723288
$ws.Range("D24").Value = 8.755368715856058
$ws.Range("E24").Value = 13.67060068609849
$ws.Range("F24").Value = 36.05169951017702
$ws.Range("J24").Value = 10.26548740095901
$ws.Range("K24").Value = 10.10684589384427
$ws.Range("M24").Value = 15.248755612328
$ws.Range("O24").Value = 27.41037948831555

# Row 25
$ws.Range("B25").Value = 10.5286337241444
$ws.Range("C25").Value = 4.546364004596921
$ws.Range("D25").Value = 8.681783727215086
$ws.Range("E25").Value = 13.64454513077255
$ws.Range("F25").Value = 36.13004440825661
$ws.Range("J25").Value = 10.29730255044692
$ws.Range("K25").Value = 9.795640552352932
$ws.Range("M25").Value = 15.1168767898513
$ws.Range("O25").Value = 27.53076307467204
